$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 25.37000000000053
$ws.Range("H2").Value = [double]"1.193788198521674e-16"
$ws.Range("K2").Value = 57.68819121972619
$ws.Range("L2").Value = "[51.81034810587181, 63.56603433358056]"
$ws.Range("P2").Value = "[1.4654476241970409, 1.6918687163476553]"
$ws.Range("S2").Value = 54.82093795740162
$ws.Range("T2").Value = "[50.76697587121813, 58.87490004358511]"
$ws.Range("W2").Value = 18.99575575575615
$ws.Range("X2").Value = 18.53863863863902
$ws.Range("Y2").Value = 19.45287287287328

# Row 3
$ws.Range("E3").Value = 22.94000000000015
$ws.Range("H3").Value = [double]"1.193788198521674e-16"
$ws.Range("K3").Value = 58.98821269634977
$ws.Range("L3").Value = "[50.20080893978297, 67.77561645291657]"
$ws.Range("O3").Value = -1.22015810770054
$ws.Range("P3").Value = "[-1.371105502467616, -1.0692107129334634]"
$ws.Range("S3").Value = 54.103334653571
$ws.Range("T3").Value = "[49.232083445295864, 58.97458586184613]"
$ws.Range("W3").Value = 4.454814814814842
$ws.Range("X3").Value = 3.903703703703731
$ws.Range("Y3").Value = 5.005925925925954
